$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meals")

$ws.Range("A21").Value = "PWXX"
$ws.Range("B21").Value = "Polish Wafer"
$ws.Range("D21").Value = 49
$ws.Range("E21").Value = 98
$ws.Range("E23").Value = 499
